$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the per-environment server detail rows (Public ip / rootpw table)
$ws.Range("A7:H8").Clear()

# Update the saved selection to H16 (last selected cell)
$ws.Range("H16").Select()
